# Applies numeric corrections to the profit-calculation columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL and WVR sheets, per the scheduled
# data-refresh run (updated market prices change the derived columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 93.31579000000001
$ws.Range("I2").Value = 95.05882
$ws.Range("J2").Value = 78.5
$ws.Range("K2").Value = 95.05882
$ws.Range("L2").Value = 78.5
$ws.Range("M2").Value = 17.94118
$ws.Range("N2").Value = -304.5
$ws.Range("H5").Value = 68.90909000000001
$ws.Range("I5").Value = 65.8
$ws.Range("K5").Value = 65.8
$ws.Range("M5").Value = 49.2
$ws.Range("H132").Value = 1961
$ws.Range("I132").Value = 1765.4
$ws.Range("K132").Value = 5296.200000000001
$ws.Range("M132").Value = -2766.200000000001
$ws.Range("H138").Value = 4031.3333
$ws.Range("I138").Value = 2598.5
$ws.Range("J138").Value = 4747.75
$ws.Range("K138").Value = 7795.5
$ws.Range("L138").Value = 14243.25
$ws.Range("M138").Value = -2655.5
$ws.Range("N138").Value = -24523.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 7745017
$ws.Range("J43").Value = 490034.5
$ws.Range("L43").Value = 490034.5
$ws.Range("N43").Value = -490660.5
$ws.Range("H44").Value = 12619.588
$ws.Range("J44").Value = 12619.588
$ws.Range("L44").Value = 12619.588
$ws.Range("N44").Value = -13595.588
$ws.Range("H61").Value = 6033.3335
$ws.Range("I61").Value = 4842.857
$ws.Range("K61").Value = 4842.857
$ws.Range("M61").Value = -4630.857
$ws.Range("H74").Value = 1810.7778
$ws.Range("I74").Value = 1810.7778
$ws.Range("K74").Value = 1810.7778
$ws.Range("M74").Value = -936.7778000000001
$ws.Range("H77").Value = 1810.7778
$ws.Range("I77").Value = 1810.7778
$ws.Range("K77").Value = 9053.889000000001
$ws.Range("M77").Value = -4685.889000000001
$ws.Range("H109").Value = 10000
$ws.Range("J109").Value = 10000
$ws.Range("L109").Value = 10000
$ws.Range("N109").Value = -12774
$ws.Range("H132").Value = 2857.375
$ws.Range("I132").Value = 2902.0667
$ws.Range("J132").Value = 2187
$ws.Range("K132").Value = 8706.2001
$ws.Range("L132").Value = 6561
$ws.Range("M132").Value = -6176.2001
$ws.Range("N132").Value = -11621
$ws.Range("H136").Value = 6033.3335
$ws.Range("I136").Value = 4842.857
$ws.Range("K136").Value = 14528.571
$ws.Range("M136").Value = -11978.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 651.3077
$ws.Range("I80").Value = 393.8889
$ws.Range("K80").Value = 393.8889
$ws.Range("M80").Value = 604.1111000000001
$ws.Range("H83").Value = 651.3077
$ws.Range("I83").Value = 393.8889
$ws.Range("K83").Value = 1969.4445
$ws.Range("M83").Value = 3022.5555
$ws.Range("H134").Value = 3701.4167
$ws.Range("I134").Value = 3701.4167
$ws.Range("K134").Value = 11104.2501
$ws.Range("M134").Value = -8569.250100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 599.8
$ws.Range("I13").Value = 399.5
$ws.Range("J13").Value = 733.3333
$ws.Range("K13").Value = 399.5
$ws.Range("L13").Value = 733.3333
$ws.Range("M13").Value = -260.5
$ws.Range("N13").Value = -1011.3333
$ws.Range("H69").Value = 32000
$ws.Range("I69").Value = 24000
$ws.Range("K69").Value = 24000
$ws.Range("M69").Value = -23251
$ws.Range("H72").Value = 32000
$ws.Range("I72").Value = 24000
$ws.Range("K72").Value = 72000
$ws.Range("M72").Value = -68256
$ws.Range("H99").Value = 2198.8
$ws.Range("I99").Value = 2123.5
$ws.Range("K99").Value = 2123.5
$ws.Range("M99").Value = -625.5
$ws.Range("H126").Value = 2198.8
$ws.Range("I126").Value = 2123.5
$ws.Range("K126").Value = 6370.5
$ws.Range("M126").Value = -3900.5
$ws.Range("H132").Value = 1649.2
$ws.Range("J132").Value = 2996
$ws.Range("L132").Value = 8988
$ws.Range("N132").Value = -14048

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 66
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 66
$ws.Range("K2").Value = 0
$ws.Range("L2").ClearContents()
$ws.Range("M2").Value = 396
$ws.Range("N2").Value = -622
$ws.Range("H5").Value = 538.6667
$ws.Range("I5").Value = 515
$ws.Range("J5").Value = 574.1667
$ws.Range("K5").Value = 1545
$ws.Range("L5").Value = 1722.5001
$ws.Range("M5").Value = -1433
$ws.Range("N5").Value = -1946.5001
$ws.Range("H17").Value = 111.75
$ws.Range("I17").Value = 31.666666
$ws.Range("J17").Value = 159.8
$ws.Range("K17").Value = 94.99999800000001
$ws.Range("L17").Value = 479.4
$ws.Range("M17").Value = 74.00000199999999
$ws.Range("N17").Value = -817.4000000000001
$ws.Range("H34").Value = 1051.5714
$ws.Range("I34").Value = 164.4
$ws.Range("J34").Value = 1544.4445
$ws.Range("K34").Value = 493.2
$ws.Range("L34").Value = 4633.333500000001
$ws.Range("M34").Value = -409.2
$ws.Range("N34").Value = -4801.333500000001
$ws.Range("H39").Value = 3462
$ws.Range("I39").Value = 310
$ws.Range("J39").Value = 4189.385
$ws.Range("K39").Value = 930
$ws.Range("L39").Value = 12568.155
$ws.Range("M39").Value = -636
$ws.Range("N39").Value = -13156.155
$ws.Range("H55").Value = 6019.9165
$ws.Range("J55").Value = 7033.9
$ws.Range("L55").Value = 21101.7
$ws.Range("N55").Value = -21455.7
$ws.Range("H75").Value = 3000
$ws.Range("J75").Value = 3000
$ws.Range("L75").Value = 9000
$ws.Range("N75").Value = -10996
$ws.Range("H78").Value = 3000
$ws.Range("J78").Value = 3000
$ws.Range("L78").Value = 27000
$ws.Range("N78").Value = -36984
$ws.Range("H135").Value = 538.6667
$ws.Range("I135").Value = 515
$ws.Range("J135").Value = 574.1667
$ws.Range("K135").Value = 4635
$ws.Range("L135").Value = 5167.5003
$ws.Range("M135").Value = -2100
$ws.Range("N135").Value = -10237.5003
$ws.Range("H136").Value = 7559.8335
$ws.Range("I136").Value = 6984
$ws.Range("K136").Value = 20952
$ws.Range("M136").Value = -15852

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1742.7142
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H52").Value = 39989
$ws.Range("I52").Value = 39989
$ws.Range("K52").Value = 39989
$ws.Range("M52").Value = -39763
$ws.Range("H54").Value = 99999
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H136").Value = 3698.353
$ws.Range("I136").Value = 2003.6
$ws.Range("K136").Value = 6010.799999999999
$ws.Range("M136").Value = -3460.799999999999
